# AutomatedTests-ComplianceReports.xlsx update
# Commit: "Added/Enabled TC207 TC238 in CRPT10"
#
# - Adds new sheet "CRPT10" (TC207 -> A, TC238 -> no status yet) after CRPT_8
# - Adds new sheet "AssetBox" (TC2197/TC2198/TC2199/TC2200/TC927 -> A) after CRPT10
# - Marks CRPT-2 TC203 as Automated (new row)
# - Marks CRPT_5 TC1524 as Fixed/Automated
# - Marks CRPT_8 TC236 as Automated (comment cleared)
# - AssetBox ends up the active/selected sheet

$wb = $excel.ActiveWorkbook

# Use the first sheet's header row as the style template for the two new
# sheets' header rows (bold/colored "TestCase | Status | Comment" style).
$headerTemplate = $wb.Worksheets.Item(1).Range("A1:C1")

# ---------------------------------------------------------------------
# 1) CRPT_5 ("CRPT_5" tab) - TC1524 row: mark fixed + automated
# ---------------------------------------------------------------------
$crpt5 = $wb.Worksheets.Item("CRPT_5")
$crpt5.Activate()
$crpt5.Range("B14").Value = "A"
$crpt5.Range("C14").Value = "Fixed - Required: CusSupervisor with assets"
$crpt5.Range("C28").Select()

# ---------------------------------------------------------------------
# 2) New sheet "CRPT10", placed right after "CRPT_8"
# ---------------------------------------------------------------------
$crpt8 = $wb.Worksheets.Item("CRPT_8")
$crpt10 = $wb.Worksheets.Add($null, $crpt8)
$crpt10.Name = "CRPT10"

$headerTemplate.Copy()
$crpt10.Range("A1:C1").PasteSpecial(-4122)
$crpt10.Range("A1").Value = "TestCase"
$crpt10.Range("B1").Value = "Status"
$crpt10.Range("C1").Value = "Comment"

$crpt10.Range("A2").Value = "TC207"
$crpt10.Range("B2").Value = "A"

$crpt10.Range("A3").Value = "TC238"

$crpt10.Range("C8").Select()

# ---------------------------------------------------------------------
# 3) New sheet "AssetBox", placed right after "CRPT10"
# ---------------------------------------------------------------------
$assetBox = $wb.Worksheets.Add($null, $crpt10)
$assetBox.Name = "AssetBox"

$headerTemplate.Copy()
$assetBox.Range("A1:C1").PasteSpecial(-4122)
$assetBox.Range("A1").Value = "TestCase"
$assetBox.Range("B1").Value = "Status"
$assetBox.Range("C1").Value = "Comment"

$assetBox.Range("A2").Value = "TC2197"
$assetBox.Range("B2").Value = "A"

$assetBox.Range("A3").Value = "TC2198"
$assetBox.Range("B3").Value = "A"

$assetBox.Range("A4").Value = "TC2199"
$assetBox.Range("B4").Value = "A"

$assetBox.Range("A5").Value = "TC2200"
$assetBox.Range("B5").Value = "A"

$assetBox.Range("A6").Value = "TC927"
$assetBox.Range("B6").Value = "A"

# ---------------------------------------------------------------------
# 4) CRPT-2 ("CRPT-2" tab) - add TC203 as a new automated row (26)
# ---------------------------------------------------------------------
$crpt2 = $wb.Worksheets.Item("CRPT-2")
$crpt2.Activate()
$crpt2.Range("A26").Value = "TC203"
$crpt2.Range("B26").Value = "A"
$crpt2.Range("B26").Select()

# ---------------------------------------------------------------------
# 5) CRPT_8 ("CRPT_8" tab) - TC236 row: mark automated, clear comment
# ---------------------------------------------------------------------
$crpt8.Activate()
$crpt8.Range("B13").Value = "A"
$crpt8.Range("C13").ClearContents()
$crpt8.Range("A1:C2").Select()

# ---------------------------------------------------------------------
# Leave AssetBox as the final active / selected sheet
# ---------------------------------------------------------------------
$assetBox.Activate()
$assetBox.Range("B12").Select()
